# Append 15 new rows (A:C) to both worksheets, matching the style of the
# existing data row (row 2), and extend the used dimension accordingly.

$wb = $excel.ActiveWorkbook

# column A | column B | column C, for each new row (rows 3..17)
$nbrData = @(
    @(1,  6,  821),
    @(2,  7,  820),
    @(3,  8,  823),
    @(4,  9,  820),
    @(5,  10, 819),
    @(6,  11, 840),
    @(7,  12, 843),
    @(8,  13, 840),
    @(9,  14, 823),
    @(10, 15, 829),
    @(11, 16, 823),
    @(12, 17, 827),
    @(13, 18, 828),
    @(14, 19, 804),
    @(15, 20, 804)
)

$barData = @(
    @(1,  6,  1207),
    @(2,  7,  1207),
    @(3,  8,  1207),
    @(4,  9,  1211),
    @(5,  10, 1209),
    @(6,  11, 1187),
    @(7,  12, 1177),
    @(8,  13, 1178),
    @(9,  14, 1169),
    @(10, 15, 1170),
    @(11, 16, 1168),
    @(12, 17, 1155),
    @(13, 18, 1148),
    @(14, 19, 1147),
    @(15, 20, 1145)
)

function Add-Rows($ws, $rows) {
    $startRow = 3
    $r = $startRow
    foreach ($row in $rows) {
        $ws.Cells.Item($r, 1).Value = $row[0]
        $ws.Cells.Item($r, 2).Value = $row[1]
        $ws.Cells.Item($r, 3).Value = $row[2]
        $r++
    }
    $lastRow = $r - 1

    # Column A on the existing data row (row 2) carries the bordered/bold
    # "s=1" style - replicate it down the new column-A cells only.
    $ws.Range("A2").Copy() | Out-Null
    $ws.Range("A" + $startRow + ":A" + $lastRow).PasteSpecial(-4122) | Out-Null
}

$wsNbr = $wb.Worksheets.Item("NBR")
Add-Rows $wsNbr $nbrData

$wsBar = $wb.Worksheets.Item("BAR")
Add-Rows $wsBar $barData

$excel.CutCopyMode = 0
